$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Cell, [string]$Text)
    $escaped = $Text -replace '"', '""'
    $Cell.Formula = '="' + $escaped + '"'
    $Cell.Copy()
    $Cell.PasteSpecial(-4163)
    $excel.CutCopyMode = $false
}

# Row 2: NCT02599584
Set-TextValue $ws.Cells.Item(2, 1) '4'
$ws.Cells.Item(2, 2).Value = '4 : pas de résultats postés ni publiés'
$ws.Cells.Item(2, 3).Value = 'NCT02599584'
Set-TextValue $ws.Cells.Item(2, 6) '2016'
$ws.Cells.Item(2, 7).Value = 'Improvement of Team Performance With a Pre-critical Staff Within the Team Prior to Admission of Critical Simulated Patient in High Fidelity Simulation Session for Anesthesiology Resident During Their Education Curriculum'
$ws.Cells.Item(2, 8).Value = 'Simstaf&perf'
$ws.Cells.Item(2, 9).Value = 'BEHAVIORAL'

# Row 3: NCT02470195
Set-TextValue $ws.Cells.Item(3, 1) '3'
$ws.Cells.Item(3, 2).Value = '3 : résultats postés ou publiés après les 36 mois'
$ws.Cells.Item(3, 3).Value = 'NCT02470195'
Set-TextValue $ws.Cells.Item(3, 6) '2016'
$ws.Cells.Item(3, 7).Value = 'Anesthesiology Education Implementation by Procedural Simulation Workshop for Difficult Airway Management: a Controlled Interregional French Study'
$ws.Cells.Item(3, 8).Value = 'PROSIDIAIR'
$ws.Cells.Item(3, 9).Value = 'BEHAVIORAL'

# Row 4: NCT02470130
Set-TextValue $ws.Cells.Item(4, 1) '4'
$ws.Cells.Item(4, 2).Value = '4 : pas de résultats postés ni publiés'
$ws.Cells.Item(4, 3).Value = 'NCT02470130'
Set-TextValue $ws.Cells.Item(4, 6) '2016'
$ws.Cells.Item(4, 7).Value = 'Effect of a Debriefing Conversational Relaxation on Memorization of Critical Key Points in Resident High Fidelity Simulation Education Program'
$ws.Cells.Item(4, 8).Value = 'RELAXSIMHF'
$ws.Cells.Item(4, 9).Value = 'BEHAVIORAL'

# Row 5: NCT02926599
Set-TextValue $ws.Cells.Item(5, 1) '4'
$ws.Cells.Item(5, 2).Value = '4 : pas de résultats postés ni publiés'
$ws.Cells.Item(5, 3).Value = 'NCT02926599'
Set-TextValue $ws.Cells.Item(5, 6) '2017'
$ws.Cells.Item(5, 7).Value = 'High Fidelity Simulation Performance After a Potential Optimization Training for Anesthesiologist Resident: a Randomized Controlled Trial.'
$ws.Cells.Item(5, 8).Value = 'TOP'
$ws.Cells.Item(5, 9).Value = 'BEHAVIORAL'

# Row 6: NCT03253770
Set-TextValue $ws.Cells.Item(6, 1) '4'
$ws.Cells.Item(6, 2).Value = '4 : pas de résultats postés ni publiés'
$ws.Cells.Item(6, 3).Value = 'NCT03253770'
Set-TextValue $ws.Cells.Item(6, 6) '2017'
$ws.Cells.Item(6, 7).Value = 'Use of a Hand-held Digital Cognitive Aid in Simulated Cardiac Arrest.'
$ws.Cells.Item(6, 8).Value = 'SIMMAX2'
$ws.Cells.Item(6, 9).Value = 'DEVICE'

# Row 7: NCT04352959
Set-TextValue $ws.Cells.Item(7, 1) '1'
$ws.Cells.Item(7, 2).Value = '1 : résultats postés ou publiés dans les 12 mois'
$ws.Cells.Item(7, 3).Value = 'NCT04352959'
Set-TextValue $ws.Cells.Item(7, 6) '2020'
$ws.Cells.Item(7, 7).Value = 'COVID-19: Nasal and Salivary Detection of the SARS-CoV-2 Virus After Antiviral Mouthrinses: Double-blind, Randomized, Placebo-controlled Clinical Study'
$ws.Cells.Item(7, 8).Value = 'BBCovid'
$ws.Cells.Item(7, 9).Value = 'DEVICE'

# Row 8: NCT04141124
Set-TextValue $ws.Cells.Item(8, 1) '2'
$ws.Cells.Item(8, 2).Value = '2 : résultats postés ou publiés entre 12 et 36 mois'
$ws.Cells.Item(8, 3).Value = 'NCT04141124'
Set-TextValue $ws.Cells.Item(8, 6) '2020'
$ws.Cells.Item(8, 7).Value = 'Effects of Relaxing Breathing Combined With Biofeedback on the Performance and Stress of Residents During a High-fidelity Simulation Session.'
$ws.Cells.Item(8, 8).Value = 'RETROSIMU'
$ws.Cells.Item(8, 9).Value = 'OTHER'

# Row 9: NCT05390879
Set-TextValue $ws.Cells.Item(9, 1) '4'
$ws.Cells.Item(9, 2).Value = '4 : pas de résultats postés ni publiés'
$ws.Cells.Item(9, 3).Value = 'NCT05390879'
Set-TextValue $ws.Cells.Item(9, 6) '2022'
$ws.Cells.Item(9, 7).Value = 'Influence of Meditation on Stress and Rumination Following Objective Structured Clinical Examination (OSCE)'
$ws.Cells.Item(9, 8).ClearContents()
$ws.Cells.Item(9, 9).Value = 'OTHER'

# Row 10: NCT05136586
Set-TextValue $ws.Cells.Item(10, 1) '4'
$ws.Cells.Item(10, 2).Value = '4 : pas de résultats postés ni publiés'
$ws.Cells.Item(10, 3).Value = 'NCT05136586'
Set-TextValue $ws.Cells.Item(10, 6) '2022'
$ws.Cells.Item(10, 7).Value = 'Effects of Two Stress Management Procedures on Performances During Objective Structured Clinical Examination (OSCE) for Medical Students : Relaxing Breathing Combined With Biofeedback or Meditative Stimulation : ECOSTRESS Study'
$ws.Cells.Item(10, 8).Value = 'ECOSTRESS'
$ws.Cells.Item(10, 9).Value = 'OTHER'

# Row 11: NCT05628519
Set-TextValue $ws.Cells.Item(11, 1) '4'
$ws.Cells.Item(11, 2).Value = '4 : pas de résultats postés ni publiés'
$ws.Cells.Item(11, 3).Value = 'NCT05628519'
Set-TextValue $ws.Cells.Item(11, 6) '2022'
$ws.Cells.Item(11, 7).Value = 'Captain Sonar Impact on Trauma Patient Management'
$ws.Cells.Item(11, 8).Value = 'CAST2'
$ws.Cells.Item(11, 9).Value = 'BEHAVIORAL'

# Row 12: NCT05619081
Set-TextValue $ws.Cells.Item(12, 1) '4'
$ws.Cells.Item(12, 2).Value = '4 : pas de résultats postés ni publiés'
$ws.Cells.Item(12, 3).Value = 'NCT05619081'
Set-TextValue $ws.Cells.Item(12, 6) '2023'
$ws.Cells.Item(12, 7).Value = 'Recovery Napping Protocol for Anesthesiologist Performance'
$ws.Cells.Item(12, 8).Value = 'R-NAP'
$ws.Cells.Item(12, 9).Value = 'BEHAVIORAL'

# Row 13: NCT05380076
Set-TextValue $ws.Cells.Item(13, 1) '4'
$ws.Cells.Item(13, 2).Value = '4 : pas de résultats postés ni publiés'
$ws.Cells.Item(13, 3).Value = 'NCT05380076'
Set-TextValue $ws.Cells.Item(13, 6) '2023'
$ws.Cells.Item(13, 7).Value = 'Effects of Mindfulness, Mobilization of Inner Resources, or Cardiac Biofeedback on Psychophysiological Anticipatory Stress Before OSCE of Medical Students'
$ws.Cells.Item(13, 8).ClearContents()
$ws.Cells.Item(13, 9).Value = 'OTHER'

# Row 14: NCT05393219
Set-TextValue $ws.Cells.Item(14, 1) '4'
$ws.Cells.Item(14, 2).Value = '4 : pas de résultats postés ni publiés'
$ws.Cells.Item(14, 3).Value = 'NCT05393219'
Set-TextValue $ws.Cells.Item(14, 6) '2023'
$ws.Cells.Item(14, 7).Value = 'Effects of Preventive Physiological and Psychological Interventions on Performances During Objective Structured Clinical Examination (OSCE) for Medical Students: Cardiac Biofeedback, Mindfulness, or Inner Resources Mobilization'
$ws.Cells.Item(14, 8).ClearContents()
$ws.Cells.Item(14, 9).Value = 'OTHER'

# Row 15: NCT06487208
Set-TextValue $ws.Cells.Item(15, 1) '4'
$ws.Cells.Item(15, 2).Value = '4 : pas de résultats postés ni publiés'
$ws.Cells.Item(15, 3).Value = 'NCT06487208'
Set-TextValue $ws.Cells.Item(15, 6) '2024'
$ws.Cells.Item(15, 7).Value = 'AnticipaMax: Impact of the Use of a Hand-held Digital Cognitive Aid in Order to Anticipated a Potential Crisis Situation, Measured by a Composite Globale Performance Score, in Anesthesia Resident'
$ws.Cells.Item(15, 8).Value = 'AnticipaMax'
$ws.Cells.Item(15, 9).Value = 'OTHER'
